$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = -11.39489999999999
$ws.Range("B12").Value = 5.774399999999999
$ws.Range("C12").Value = -14.69310000000002
$ws.Range("C14").Value = -12.10289999999999
$ws.Range("C22").Value = -10.81769999999999
